# Final Report Statistics - figure refresh
#
# Updates the "Realized" figures that feed the two charts on Sheet1
# (Monthly Planned vs. Realized Time Usage, and Working Hour
# Distribution) plus the underlying Table1 data, then leaves the
# selection where the author left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5 ("Realized") feeds Chart 1 (Sheet1!$B$5:$G$5) - Mar/Apr 2013 figures
$ws.Range("F5").Value = 93.55
$ws.Range("G5").Value = 55.5

# Column M ("Realized") feeds Chart 2 / Table1 (Sheet1!$M$4:$M$10)
$ws.Range("M4").Value = 96.3
$ws.Range("M5").Value = 41
$ws.Range("M7").Value = 88.25

# Leave the selection on M5, matching where editing finished
$ws.Range("M5").Select()
